$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'323.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'8.78%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'49.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'18.98%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.352"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'7.04%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08168"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'8.62%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.608"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'5.29%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.668"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'5.30%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D9").Value = "'0.1354"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'13.41%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1955"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'7.25%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09515"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'6.55%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04558"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'11.79%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1049"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.05%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001319"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'2.81%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005946"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'2.74%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.399"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.25%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D18").Value = "'0.3393"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.41%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'8.216"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.40%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1418"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'4.17%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.3053"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.60%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04302"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'5.20%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001305"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'3.06%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004268"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'9.36%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001350"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'9.74%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003722"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.06%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02778"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'15.20%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05559"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'6.64%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.006299"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.09%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007710"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.58%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1450"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'9.38%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007693"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'3.88%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008082"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'11.29%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3512"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'18.34%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006771"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'2.99%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.09%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D49").Value = "'0.004000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-4.84%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.09%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.09%"
$ws.Range("E51").Style = "Normal"
